$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#mena"
$ws.Range("C2").Value = "Mena"
$ws.Range("D2").ClearContents()
$ws.Range("B3").Value = "#ulyssis"
$ws.Range("C3").Value = "Ulyssis"
$ws.Range("D3").ClearContents()
$ws.Range("B4").Value = "#eli"
$ws.Range("C4").Value = "Eli"
$ws.Range("D4").ClearContents()
$ws.Range("B5").Value = "#dienaer"
$ws.Range("C5").Value = "Dienaer"
$ws.Range("D5").ClearContents()
$ws.Range("B6").Value = "#pria"
$ws.Range("C6").Value = "Pria"
$ws.Range("D6").ClearContents()
$ws.Range("B7").Value = "#ayax"
$ws.Range("C7").Value = "Ayax"
$ws.Range("D7").ClearContents()
$ws.Range("B8").Value = "#menand"
$ws.Range("C8").Value = "Menand"
$ws.Range("D8").ClearContents()
$ws.Range("B9").Value = "#priamus"
$ws.Range("C9").Value = "Priamus"
$ws.Range("D9").ClearContents()
$ws.Range("B10").Value = "#page"
$ws.Range("C10").Value = "Page"
$ws.Range("D10").ClearContents()
$ws.Range("B11").Value = "#edelen"
$ws.Range("C11").Value = "Edelen"
$ws.Range("D11").ClearContents()
$ws.Range("B12").Value = "#venus-lachende"
$ws.Range("C12").Value = "Venus lachende"
$ws.Range("D12").ClearContents()
$ws.Range("B13").Value = "#menelaus"
$ws.Range("C13").Value = "Menelaus"
$ws.Range("D13").ClearContents()
$ws.Range("B14").Value = "#rey"
$ws.Range("C14").Value = "Rey"
$ws.Range("D14").ClearContents()
$ws.Range("B15").Value = "#venus"
$ws.Range("C15").Value = "Venus"
$ws.Range("D15").ClearContents()
$ws.Range("B16").Value = "#iuno"
$ws.Range("C16").Value = "Iuno"
$ws.Range("D16").ClearContents()
$ws.Range("B17").Value = "#ajax"
$ws.Range("C17").Value = "Ajax"
$ws.Range("D17").ClearContents()
$ws.Range("B18").Value = "#hele"
$ws.Range("C18").Value = "Hele"
$ws.Range("D18").ClearContents()
$ws.Range("B19").Value = "#oeno"
$ws.Range("C19").Value = "OEno"
$ws.Range("D19").ClearContents()
$ws.Range("B20").Value = "#pallas"
$ws.Range("C20").Value = "Pallas"
$ws.Range("D20").ClearContents()
$ws.Range("B21").Value = "#oenone"
$ws.Range("C21").Value = "OEnone"
$ws.Range("D21").ClearContents()
$ws.Range("B22").Value = "#mercur"
$ws.Range("C22").Value = "Mercur"
$ws.Range("D22").ClearContents()
$ws.Range("B23").Value = "#helena"
$ws.Range("C23").Value = "Helena"
$ws.Range("D23").ClearContents()
$ws.Range("B24").Value = "#pagie"
$ws.Range("C24").Value = "Pagie"
$ws.Range("D24").ClearContents()
$ws.Range("B25").Value = "#ulys"
$ws.Range("C25").Value = "Ulys"
$ws.Range("D25").ClearContents()
$ws.Range("B26").Value = "#elicia"
$ws.Range("C26").Value = "Elicia"
$ws.Range("D26").ClearContents()
$ws.Range("B27").Value = "#mene"
$ws.Range("C27").Value = "Mene"
$ws.Range("D27").ClearContents()
$ws.Range("B28").Value = "#paris"
$ws.Range("C28").Value = "Paris"
$ws.Range("D28").ClearContents()
$ws.Range("B29").Value = "#par"
$ws.Range("C29").Value = "Par"
$ws.Range("D29").ClearContents()
